# Auto-generated edit script applying the Bahamut_Profits.xlsx leve price/profit refresh
# (currentAveragePrice* / LevePrice* / LeveProfit* columns) across the ALC, ARM, CRP, CUL,
# GSM, LTW and WVR sheets, per the scheduled price-refresh runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 10000
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 10000
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H87").Value = 38000
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 38000
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 38000
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -40496
$ws.Range("H90").Value = 38000
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 38000
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 114000
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -126480
$ws.Range("H116").Value = 3971.08
$ws.Range("I116").Value = 3973.9
$ws.Range("K116").Value = 3973.9
$ws.Range("M116").Value = -531.9000000000001
$ws.Range("H129").Value = 772637.4
$ws.Range("I129").Value = 431
$ws.Range("J129").Value = 904477.4399999999
$ws.Range("K129").Value = 1293
$ws.Range("L129").Value = 2713432.32
$ws.Range("M129").Value = 3707
$ws.Range("N129").Value = -2723432.32
$ws.Range("H137").Value = 826.4048
$ws.Range("I137").Value = 749.7353000000001
$ws.Range("J137").Value = 1152.25
$ws.Range("K137").Value = 2249.2059
$ws.Range("L137").Value = 3456.75
$ws.Range("M137").Value = 300.7941000000001
$ws.Range("N137").Value = -8556.75
$ws.Range("H138").Value = 2480.78
$ws.Range("I138").Value = 1080.6
$ws.Range("J138").Value = 3414.2334
$ws.Range("K138").Value = 3241.8
$ws.Range("L138").Value = 10242.7002
$ws.Range("M138").Value = 1898.2
$ws.Range("N138").Value = -20522.7002
$ws.Range("H139").Value = 94950
$ws.Range("J139").Value = 94950
$ws.Range("L139").Value = 94950
$ws.Range("N139").Value = -105230
$ws.Range("H140").Value = 51617.777
$ws.Range("J140").Value = 51617.777
$ws.Range("L140").Value = 51617.777
$ws.Range("N140").Value = -61977.777

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6285.66
$ws.Range("I32").Value = 6160.4893
$ws.Range("J32").Value = 8246.666999999999
$ws.Range("K32").Value = 6160.4893
$ws.Range("L32").Value = 8246.666999999999
$ws.Range("M32").Value = -5873.4893
$ws.Range("N32").Value = -8820.666999999999
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H139").Value = 68919.164
$ws.Range("J139").Value = 68919.164
$ws.Range("L139").Value = 68919.164
$ws.Range("N139").Value = -79199.164
$ws.Range("H140").Value = 57124
$ws.Range("J140").Value = 57124
$ws.Range("L140").Value = 57124
$ws.Range("N140").Value = -67484
$ws.Range("H141").Value = 47073.715
$ws.Range("J141").Value = 47073.715
$ws.Range("L141").Value = 47073.715
$ws.Range("N141").Value = -57433.715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1539.125
$ws.Range("I22").Value = 1801.8334
$ws.Range("J22").Value = 751
$ws.Range("K22").Value = 1801.8334
$ws.Range("L22").Value = 751
$ws.Range("M22").Value = -1451.8334
$ws.Range("N22").Value = -1451
$ws.Range("H31").Value = 86936.914
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 86936.914
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 86936.914
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -87526.914
$ws.Range("H34").Value = 86936.914
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 86936.914
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 86936.914
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -87340.914
$ws.Range("H41").Value = 14599.2
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 14599.2
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 14599.2
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -15455.2
$ws.Range("H50").Value = 9232.6
$ws.Range("J50").Value = 9232.6
$ws.Range("L50").Value = 9232.6
$ws.Range("N50").Value = -10482.6
$ws.Range("H51").Value = 8517.556
$ws.Range("J51").Value = 8866.333000000001
$ws.Range("L51").Value = 8866.333000000001
$ws.Range("N51").Value = -10338.333
$ws.Range("H58").Value = 3358.1162
$ws.Range("I58").Value = 917.6786
$ws.Range("J58").Value = 7913.6
$ws.Range("K58").Value = 917.6786
$ws.Range("L58").Value = 7913.6
$ws.Range("M58").Value = -714.6786
$ws.Range("N58").Value = -8319.6
$ws.Range("H59").Value = 16010.8
$ws.Range("J59").Value = 16010.8
$ws.Range("L59").Value = 16010.8
$ws.Range("N59").Value = -18300.8
$ws.Range("H60").Value = 7251.3335
$ws.Range("J60").Value = 8101.6
$ws.Range("L60").Value = 8101.6
$ws.Range("N60").Value = -9123.6
$ws.Range("H61").Value = 8517.556
$ws.Range("J61").Value = 8866.333000000001
$ws.Range("L61").Value = 8866.333000000001
$ws.Range("N61").Value = -9562.333000000001
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496
$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716
$ws.Range("H130").Value = 45000
$ws.Range("J130").Value = 45000
$ws.Range("L130").Value = 45000
$ws.Range("N130").Value = -55040
$ws.Range("H136").Value = 3358.1162
$ws.Range("I136").Value = 917.6786
$ws.Range("J136").Value = 7913.6
$ws.Range("K136").Value = 2753.0358
$ws.Range("L136").Value = 23740.8
$ws.Range("M136").Value = -203.0357999999997
$ws.Range("N136").Value = -28840.8
$ws.Range("H138").Value = 77882.22
$ws.Range("J138").Value = 77882.22
$ws.Range("L138").Value = 77882.22
$ws.Range("N138").Value = -88162.22
$ws.Range("H140").Value = 69280
$ws.Range("J140").Value = 69280
$ws.Range("L140").Value = 69280
$ws.Range("N140").Value = -79640

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 100165.2
$ws.Range("J29").Value = 125181.5
$ws.Range("L29").Value = 375544.5
$ws.Range("N29").Value = -376098.5
$ws.Range("H69").Value = 452.16666
$ws.Range("J69").Value = 600
$ws.Range("L69").Value = 1800
$ws.Range("N69").Value = -3422
$ws.Range("H72").Value = 452.16666
$ws.Range("J72").Value = 600
$ws.Range("L72").Value = 5400
$ws.Range("N72").Value = -13512
$ws.Range("H113").Value = 622.2414
$ws.Range("I113").Value = 941.2857
$ws.Range("J113").Value = 520.7273
$ws.Range("K113").Value = 2823.8571
$ws.Range("L113").Value = 1562.1819
$ws.Range("M113").Value = -653.8571000000002
$ws.Range("N113").Value = -5902.1819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4352.125
$ws.Range("I70").Value = 4136.875
$ws.Range("J70").Value = 4567.375
$ws.Range("K70").Value = 4136.875
$ws.Range("L70").Value = 4567.375
$ws.Range("M70").Value = -3866.875
$ws.Range("N70").Value = -5107.375
$ws.Range("H73").Value = 4352.125
$ws.Range("I73").Value = 4136.875
$ws.Range("J73").Value = 4567.375
$ws.Range("K73").Value = 4136.875
$ws.Range("L73").Value = 4567.375
$ws.Range("M73").Value = -3200.875
$ws.Range("N73").Value = -6439.375
$ws.Range("H132").Value = 3064.2979
$ws.Range("I132").Value = 2818.8333
$ws.Range("J132").Value = 3867.6365
$ws.Range("K132").Value = 8456.499899999999
$ws.Range("L132").Value = 11602.9095
$ws.Range("M132").Value = -5926.499899999999
$ws.Range("N132").Value = -16662.9095
$ws.Range("H140").Value = 74220
$ws.Range("J140").Value = 74220
$ws.Range("L140").Value = 74220
$ws.Range("N140").Value = -84580
$ws.Range("H141").Value = 59933.332
$ws.Range("J141").Value = 59933.332
$ws.Range("L141").Value = 59933.332
$ws.Range("N141").Value = -70293.33199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2495.8333
$ws.Range("I7").Value = 2387.652
$ws.Range("K7").Value = 2387.652
$ws.Range("M7").Value = -2275.652
$ws.Range("H16").Value = 2034.75
$ws.Range("I16").Value = 2110.4546
$ws.Range("J16").Value = 1202
$ws.Range("K16").Value = 2110.4546
$ws.Range("L16").Value = 1202
$ws.Range("M16").Value = -1940.4546
$ws.Range("N16").Value = -1542
$ws.Range("H126").Value = 2495.8333
$ws.Range("I126").Value = 2387.652
$ws.Range("K126").Value = 7162.956
$ws.Range("M126").Value = -4692.956
$ws.Range("H136").Value = 2437.7932
$ws.Range("I136").Value = 1482.625
$ws.Range("J136").Value = 5095.6523
$ws.Range("K136").Value = 4447.875
$ws.Range("L136").Value = 15286.9569
$ws.Range("M136").Value = -1897.875
$ws.Range("N136").Value = -20386.9569
$ws.Range("H140").Value = 77994.60000000001
$ws.Range("J140").Value = 77994.60000000001
$ws.Range("L140").Value = 77994.60000000001
$ws.Range("N140").Value = -88354.60000000001
$ws.Range("H141").Value = 67138.336
$ws.Range("J141").Value = 67138.336
$ws.Range("L141").Value = 67138.336
$ws.Range("N141").Value = -77498.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8583.333000000001
$ws.Range("J62").Value = 9125
$ws.Range("L62").Value = 9125
$ws.Range("N62").Value = -10373
$ws.Range("H65").Value = 8583.333000000001
$ws.Range("J65").Value = 9125
$ws.Range("L65").Value = 45625
$ws.Range("N65").Value = -51865
$ws.Range("H136").Value = 696.75
$ws.Range("I136").Value = 586.1667
$ws.Range("J136").Value = 1249.6666
$ws.Range("K136").Value = 1758.5001
$ws.Range("L136").Value = 3748.9998
$ws.Range("M136").Value = 791.4999
$ws.Range("N136").Value = -8848.9998
$ws.Range("H140").Value = 77828.664
$ws.Range("J140").Value = 77828.664
$ws.Range("L140").Value = 77828.664
$ws.Range("N140").Value = -88188.664
$ws.Range("H141").Value = 89533.89
$ws.Range("J141").Value = 89533.89
$ws.Range("L141").Value = 89533.89
$ws.Range("N141").Value = -99893.89
